# Actualización automática 2025-07-31 17:10:08
#
# Updates a single underlying sale amount (+3858.62) for
# RIOS CARRION ANGEL BENIGNO / PORCELANATO, which cascades into the
# group-sales sheet, the monthly-sales sheet and the compliance sheet.

$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- VENTAS POR GRUPO: PORCELANATO sale for RIOS CARRION ANGEL BENIGNO ---
$wsGrupo.Range("M11").Value = 3864.56

# --- VENTA MENSUAL: julio sale for the same client, plus its total row ---
$wsMensual.Range("F11").Value = 5087.14
$wsMensual.Range("F24").Value = 54512.04

# --- CUMPLIMIENTO MENSUAL: PORCELANATO row and TOTAL row ---
$wsCumplimiento.Range("D16").Value = 49179.22
$wsCumplimiento.Range("E16").Value = -10422.68
$wsCumplimiento.Range("F16").Value = 1.268926999159368

$wsCumplimiento.Range("D19").Value = 54512.04
$wsCumplimiento.Range("E19").Value = 3710.963863046034
$wsCumplimiento.Range("F19").Value = 0.9362629267329615

# Column D on "CUMPLIMIENTO MENSUAL" widens slightly to fit the new value.
# (13.15 is used instead of 14 because the COM ColumnWidth setter applies an
# internal pixel-rounding offset; this value is the one that serializes back
# to the clean OOXML width of 14.)
$wsCumplimiento.Columns.Item(4).ColumnWidth = 13.15
